# Jogos_da_Semana_FlashScore_2024-10-13.xlsx -- apply odds updates, then
# remove the stale "Oakland Roots vs Phoenix Rising" fixture (row 24),
# which shifts the remaining USA-USL/VENEZUELA rows up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("G5").Value = 1.85
$ws.Range("H5").Value = 3.7
$ws.Range("I5").Value = 3.9
$ws.Range("J5").Value = 2.6
$ws.Range("L5").Value = 5
$ws.Range("N5").Value = 8.5
$ws.Range("Q5").Value = 2.2
$ws.Range("R5").Value = 1.65
$ws.Range("AG5").Value = 9
$ws.Range("AO5").Value = 10
$ws.Range("AQ5").Value = 34
$ws.Range("AX5").Value = 26

# Row 7
$ws.Range("G7").Value = 1.36
$ws.Range("H7").Value = 4.6
$ws.Range("I7").Value = 7.2
$ws.Range("K7").Value = 2.45
$ws.Range("L7").Value = 6.3
$ws.Range("N7").Value = 13.9
$ws.Range("Q7").Value = 1.53
$ws.Range("R7").Value = 2.18
$ws.Range("S7").Value = 1.28
$ws.Range("T7").Value = 3.46
$ws.Range("U7").Value = 1.78
$ws.Range("V7").Value = 1.83
$ws.Range("W7").Value = 7.9
$ws.Range("X7").Value = 7
$ws.Range("Z7").Value = 9
$ws.Range("AA7").Value = 10.75
$ws.Range("AB7").Value = 24
$ws.Range("AC7").Value = 14.5
$ws.Range("AD7").Value = 9.25
$ws.Range("AE7").Value = 18.5
$ws.Range("AG7").Value = 22
$ws.Range("AH7").Value = 50
$ws.Range("AK7").Value = 75
$ws.Range("AO7").Value = 6
$ws.Range("AP7").Value = 15
$ws.Range("AR7").Value = 40
$ws.Range("AU7").Value = 7.9
$ws.Range("AW7").Value = 8.25
$ws.Range("AY7").Value = 37
$ws.Range("AZ7").Value = 250
$ws.Range("BA7").Value = 250
$ws.Range("BB7").Value = 450

# Row 8
$ws.Range("T8").Value = 2.97

# Row 9
$ws.Range("G9").Value = 2.1
$ws.Range("I9").Value = 3.1
$ws.Range("L9").Value = 3.5
$ws.Range("X9").Value = 13
$ws.Range("Z9").Value = 21
$ws.Range("AC9").Value = 15
$ws.Range("AD9").Value = 7
$ws.Range("AE9").Value = 11
$ws.Range("AX9").Value = 15

# Row 12
$ws.Range("I12").Value = 2.9
$ws.Range("K12").Value = 2.05
$ws.Range("U12").Value = 1.83
$ws.Range("V12").Value = 1.83
$ws.Range("AB12").Value = 34
$ws.Range("AG12").Value = 8.5
$ws.Range("AK12").Value = 26
$ws.Range("AM12").Value = 301
$ws.Range("AS12").Value = 201

# Row 13
$ws.Range("N13").Value = 9

# Row 15
$ws.Range("G15").Value = 2.5
$ws.Range("H15").Value = 3
$ws.Range("I15").Value = 3
$ws.Range("K15").Value = 2
$ws.Range("M15").Value = 1.1
$ws.Range("N15").Value = 7
$ws.Range("S15").Value = 1.5
$ws.Range("T15").Value = 2.5
$ws.Range("AD15").Value = 6
$ws.Range("AH15").Value = 13
$ws.Range("AJ15").Value = 29
$ws.Range("AM15").Value = 900
$ws.Range("AT15").Value = 2.5
$ws.Range("AX15").Value = 17
$ws.Range("AZ15").Value = 51
$ws.Range("BA15").Value = 81

# Row 16
$ws.Range("Q16").Value = 2.08
$ws.Range("R16").Value = 1.73

# Row 17
$ws.Range("G17").Value = 1.91
$ws.Range("H17").Value = 3.6
$ws.Range("I17").Value = 3.75
$ws.Range("J17").Value = 2.5
$ws.Range("L17").Value = 4
$ws.Range("M17").Value = 1.03
$ws.Range("N17").Value = 15
$ws.Range("U17").Value = 1.62
$ws.Range("V17").Value = 2.2
$ws.Range("W17").Value = 9.5
$ws.Range("X17").Value = 11
$ws.Range("Y17").Value = 8.5
$ws.Range("Z17").Value = 17
$ws.Range("AC17").Value = 15
$ws.Range("AE17").Value = 13
$ws.Range("AH17").Value = 21
$ws.Range("AI17").Value = 13
$ws.Range("AJ17").Value = 41
$ws.Range("AK17").Value = 26
$ws.Range("AM17").Value = 151
$ws.Range("AO17").Value = 10
$ws.Range("AP17").Value = 17
$ws.Range("AR17").Value = 41
$ws.Range("AW17").Value = 6
$ws.Range("AX17").Value = 19
$ws.Range("AY17").Value = 23
$ws.Range("BB17").Value = 151

# Row 18
$ws.Range("G18").Value = 2.15
$ws.Range("I18").Value = 3.75
$ws.Range("J18").Value = 3
$ws.Range("K18").Value = 1.91
$ws.Range("X18").Value = 9
$ws.Range("AK18").Value = 34
$ws.Range("AO18").Value = 13
$ws.Range("AU18").Value = 9
$ws.Range("AX18").Value = 21

# Row 19
$ws.Range("G19").Value = 1.91
$ws.Range("H19").Value = 3.4
$ws.Range("I19").Value = 4
$ws.Range("J19").Value = 2.6
$ws.Range("L19").Value = 4.33
$ws.Range("W19").Value = 7.5
$ws.Range("Y19").Value = 9
$ws.Range("Z19").Value = 17
$ws.Range("AW19").Value = 5.5
$ws.Range("BA19").Value = 81

# Row 20
$ws.Range("AE20").Value = 15
$ws.Range("AG20").Value = 9
$ws.Range("AQ20").Value = 51
$ws.Range("AV20").Value = 67
$ws.Range("BB20").Value = 251

# Row 21
$ws.Range("O21").Value = 1.5
$ws.Range("P21").Value = 2.5
$ws.Range("Q21").Value = 2.5
$ws.Range("R21").Value = 1.5

# Delete the old row 24 entirely; rows 25-28 shift up to become 24-27
# (dimension becomes A1:BD27).
$ws.Rows(24).Delete()
